$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 84
$ws.Range("I11").Value = 84
$ws.Range("K11").Value = 84
$ws.Range("M11").Value = 56
$ws.Range("H40").Value = 2595.1667
$ws.Range("I40").Value = 1521.0667
$ws.Range("J40").Value = 7965.6665
$ws.Range("K40").Value = 1521.0667
$ws.Range("L40").Value = 7965.6665
$ws.Range("M40").Value = -1346.0667
$ws.Range("N40").Value = -8315.666499999999
$ws.Range("H64").Value = 6488.857
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 6488.857
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 6488.857
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -6984.857
$ws.Range("H67").Value = 6488.857
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 6488.857
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 6488.857
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -8204.857
$ws.Range("H88").Value = 9227.25
$ws.Range("J88").Value = 8580.888999999999
$ws.Range("L88").Value = 8580.888999999999
$ws.Range("N88").Value = -9392.888999999999
$ws.Range("H91").Value = 9227.25
$ws.Range("J91").Value = 8580.888999999999
$ws.Range("L91").Value = 8580.888999999999
$ws.Range("N91").Value = -11388.889
$ws.Range("H97").Value = 4016.3333
$ws.Range("J97").Value = 4016.3333
$ws.Range("L97").Value = 12048.9999
$ws.Range("N97").Value = -13040.9999
$ws.Range("H131").Value = 5559.25
$ws.Range("J131").Value = 20000
$ws.Range("L131").Value = 60000
$ws.Range("N131").Value = -70080
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2593.68
$ws.Range("I45").Value = 1108.4445
$ws.Range("K45").Value = 1108.4445
$ws.Range("M45").Value = -731.4445000000001
$ws.Range("H61").Value = 5696.0835
$ws.Range("I61").Value = 4917.875
$ws.Range("K61").Value = 4917.875
$ws.Range("M61").Value = -4705.875
$ws.Range("H92").Value = 98316
$ws.Range("J92").Value = 98316
$ws.Range("L92").Value = 98316
$ws.Range("N92").Value = -103308
$ws.Range("H132").Value = 45149.043
$ws.Range("I132").Value = 61919.35
$ws.Range("J132").Value = 4421.143
$ws.Range("K132").Value = 185758.05
$ws.Range("L132").Value = 13263.429
$ws.Range("M132").Value = -183228.05
$ws.Range("N132").Value = -18323.429
$ws.Range("H136").Value = 5696.0835
$ws.Range("I136").Value = 4917.875
$ws.Range("K136").Value = 14753.625
$ws.Range("M136").Value = -12203.625
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3706.8096
$ws.Range("I134").Value = 2709
$ws.Range("J134").Value = 6899.8
$ws.Range("K134").Value = 8127
$ws.Range("L134").Value = 20699.4
$ws.Range("M134").Value = -5592
$ws.Range("N134").Value = -25769.4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3817.2307
$ws.Range("I62").Value = 3494
$ws.Range("J62").Value = 4194.3335
$ws.Range("K62").Value = 3494
$ws.Range("L62").Value = 4194.3335
$ws.Range("M62").Value = -2870
$ws.Range("N62").Value = -5442.3335
$ws.Range("H65").Value = 3817.2307
$ws.Range("I65").Value = 3494
$ws.Range("J65").Value = 4194.3335
$ws.Range("K65").Value = 17470
$ws.Range("L65").Value = 20971.6675
$ws.Range("M65").Value = -14350
$ws.Range("N65").Value = -27211.6675
$ws.Range("H99").Value = 6999
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 6999
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 6999
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -9995
$ws.Range("H105").Value = 1887.591
$ws.Range("I105").Value = 1834.0952
$ws.Range("K105").Value = 1834.0952
$ws.Range("M105").Value = -87.09519999999998
$ws.Range("H126").Value = 6999
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 6999
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 20997
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -25937
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 565.35297
$ws.Range("I8").Value = 565.35297
$ws.Range("K8").Value = 1696.05891
$ws.Range("M8").Value = -1557.05891
$ws.Range("H129").Value = 503099.94
$ws.Range("I129").Value = 2565.2856
$ws.Range("K129").Value = 7695.8568
$ws.Range("M129").Value = -2695.8568
$ws.Range("H131").Value = 26018.625
$ws.Range("I131").Value = 1049.6666
$ws.Range("J131").Value = 41000
$ws.Range("K131").Value = 3148.9998
$ws.Range("L131").Value = 123000
$ws.Range("M131").Value = 1891.0002
$ws.Range("N131").Value = -133080
$ws.Range("H132").Value = 2193
$ws.Range("I132").Value = 2292.6
$ws.Range("J132").Value = 1993.8
$ws.Range("K132").Value = 20633.4
$ws.Range("L132").Value = 17944.2
$ws.Range("M132").Value = -18103.4
$ws.Range("N132").Value = -23004.2
$ws.Range("H140").Value = 5748.5
$ws.Range("I140").Value = 5498.3335
$ws.Range("K140").Value = 16495.0005
$ws.Range("M140").Value = -11315.0005
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3629.1333
$ws.Range("I80").Value = 3205.3
$ws.Range("K80").Value = 3205.3
$ws.Range("M80").Value = -2207.3
$ws.Range("H83").Value = 3629.1333
$ws.Range("I83").Value = 3205.3
$ws.Range("K83").Value = 16026.5
$ws.Range("M83").Value = -11034.5
$ws.Range("H92").Value = 26308.75
$ws.Range("J92").Value = 26308.75
$ws.Range("L92").Value = 26308.75
$ws.Range("N92").Value = -30052.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 762.1
$ws.Range("I55").Value = 467.57144
$ws.Range("J55").Value = 1449.3334
$ws.Range("K55").Value = 467.57144
$ws.Range("L55").Value = 1449.3334
$ws.Range("M55").Value = -294.57144
$ws.Range("N55").Value = -1795.3334
$ws.Range("H68").Value = 6166
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 6166
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 6166
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -7664
$ws.Range("H71").Value = 6166
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 6166
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 30830
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -38318
$ws.Range("H100").Value = 2169.611
$ws.Range("I100").Value = 1173.6154
$ws.Range("K100").Value = 1173.6154
$ws.Range("M100").Value = -632.6153999999999
$ws.Range("H136").Value = 6054.2856
$ws.Range("I136").Value = 4444
$ws.Range("J136").Value = 6698.4
$ws.Range("K136").Value = 13332
$ws.Range("L136").Value = 20095.2
$ws.Range("M136").Value = -10782
$ws.Range("N136").Value = -25195.2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2751.625
$ws.Range("J81").Value = 2222
$ws.Range("L81").Value = 4444
$ws.Range("N81").Value = -6566
$ws.Range("H84").Value = 2751.625
$ws.Range("J84").Value = 2222
$ws.Range("L84").Value = 22220
$ws.Range("N84").Value = -32828
$ws.Range("H96").Value = 3418.7273
$ws.Range("I96").Value = 2464.8
$ws.Range("J96").Value = 4213.6665
$ws.Range("K96").Value = 2464.8
$ws.Range("L96").Value = 4213.6665
$ws.Range("M96").Value = -1091.8
$ws.Range("N96").Value = -6959.6665
$ws.Range("H122").Value = 2063.7144
$ws.Range("I122").Value = 2089.4
$ws.Range("J122").Value = 1999.5
$ws.Range("K122").Value = 6268.200000000001
$ws.Range("L122").Value = 5998.5
$ws.Range("M122").Value = -3818.200000000001
$ws.Range("N122").Value = -10898.5
$ws.Range("H132").Value = 44959.832
$ws.Range("I132").Value = 46862.566
$ws.Range("J132").Value = 1197
$ws.Range("K132").Value = 140587.698
$ws.Range("L132").Value = 3591
$ws.Range("M132").Value = -138057.698
$ws.Range("N132").Value = -8651
